# B1--and-B2-PowerPoint.pptx edit
# 1) Slide 5's table switches from the custom "Table_0" style to the
#    built-in "Medium Style 2 - Accent 1" table style.
# 2) The deck's theme colour scheme (used by the slide master, theme1.xml)
#    is swapped from the "Red Violet" Integral palette to the stock
#    "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{E27873B3-6ECE-476A-ACF2-02C52EDD0601}")

# --- 2. Swap the theme colour scheme to the default "Office" palette ------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
$colorScheme.Colors(1).RGB  = 0        # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477  # folHlink 954F72
